# Scheduled-runner refresh: overwrite the market-price / profit columns
# (H:N -> currentAveragePrice.. LeveProfitHQ) with freshly pulled values
# for the rows whose quoted prices moved since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1138.4615
$arr[0,1] = 941.6667
$arr[0,2] = 3500
$arr[0,3] = 941.6667
$arr[0,4] = 3500
$arr[0,5] = 556.3333
$arr[0,6] = -6496
$ws.Range("H98:N98").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1138.4615
$arr[0,1] = 941.6667
$arr[0,2] = 3500
$arr[0,3] = 2825.0001
$arr[0,4] = 10500
$arr[0,5] = -375.0001000000002
$arr[0,6] = -15400
$ws.Range("H122:N122").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5056.625
$arr[0,1] = 1862
$arr[0,2] = 18900
$arr[0,3] = 5586
$arr[0,4] = 56700
$arr[0,5] = -3056
$arr[0,6] = -61760
$ws.Range("H132:N132").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 7221.643
$arr[0,1] = 30000
$arr[0,2] = 3425.25
$arr[0,3] = 90000
$arr[0,4] = 10275.75
$arr[0,5] = -87450
$arr[0,6] = -15375.75
$ws.Range("H137:N137").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1679.4777
$arr[0,1] = 757.5641000000001
$arr[0,2] = 2963.5715
$arr[0,3] = 2272.6923
$arr[0,4] = 8890.7145
$arr[0,5] = 2867.3077
$arr[0,6] = -19170.7145
$ws.Range("H138:N138").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1816686.2
$arr[0,1] = 1816686.2
$arr[0,2] = 0
$arr[0,3] = 1816686.2
$arr[0,4] = 0
$arr[0,5] = -1816399.2
$ws.Range("H32:M32").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1287.4166
$arr[0,1] = 1268.5
$arr[0,2] = 1306.3334
$arr[0,3] = 1268.5
$arr[0,4] = 1306.3334
$arr[0,5] = -891.5
$arr[0,6] = -2060.3334
$ws.Range("H45:N45").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4778.3687
$arr[0,1] = 1570.6428
$arr[0,2] = 13760
$arr[0,3] = 1570.6428
$arr[0,4] = 13760
$arr[0,5] = 51.35719999999992
$arr[0,6] = -17004
$ws.Range("H102:N102").Value = $arr

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1289.174
$arr[0,1] = 1101
$arr[0,2] = 1581.8889
$arr[0,3] = 1101
$arr[0,4] = 1581.8889
$arr[0,5] = -854
$arr[0,6] = -2075.8889
$ws.Range("H20:N20").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 593.7406999999999
$arr[0,1] = 522.35297
$arr[0,2] = 715.1
$arr[0,3] = 522.35297
$arr[0,4] = 715.1
$arr[0,5] = -297.35297
$arr[0,6] = -1165.1
$ws.Range("H64:N64").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 593.7406999999999
$arr[0,1] = 522.35297
$arr[0,2] = 715.1
$arr[0,3] = 522.35297
$arr[0,4] = 715.1
$arr[0,5] = 257.64703
$arr[0,6] = -2275.1
$ws.Range("H67:N67").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4535.5713
$arr[0,1] = 7828.2856
$arr[0,2] = 1242.8572
$arr[0,3] = 7828.2856
$arr[0,4] = 1242.8572
$arr[0,5] = -6330.2856
$arr[0,6] = -4238.8572
$ws.Range("H99:N99").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2779330.5
$arr[0,1] = 1383.6364
$arr[0,2] = 7144675.5
$arr[0,3] = 1383.6364
$arr[0,4] = 7144675.5
$arr[0,5] = 363.3635999999999
$arr[0,6] = -7148169.5
$ws.Range("H105:N105").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1609.9474
$arr[0,1] = 1352.091
$arr[0,2] = 1964.5
$arr[0,3] = 1352.091
$arr[0,4] = 1964.5
$arr[0,5] = 567.9090000000001
$arr[0,6] = -5804.5
$ws.Range("H107:N107").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4526.722
$arr[0,1] = 1664.5555
$arr[0,2] = 7388.8887
$arr[0,3] = 1664.5555
$arr[0,4] = 7388.8887
$arr[0,5] = -1369.5555
$arr[0,6] = -7978.8887
$ws.Range("H31:N31").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4526.722
$arr[0,1] = 1664.5555
$arr[0,2] = 7388.8887
$arr[0,3] = 1664.5555
$arr[0,4] = 7388.8887
$arr[0,5] = -1462.5555
$arr[0,6] = -7792.8887
$ws.Range("H34:N34").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2781337.8
$arr[0,1] = 6947632
$arr[0,2] = 3808.3333
$arr[0,3] = 6947632
$arr[0,4] = 3808.3333
$arr[0,5] = -6947008
$arr[0,6] = -5056.3333
$ws.Range("H62:N62").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2781337.8
$arr[0,1] = 6947632
$arr[0,2] = 3808.3333
$arr[0,3] = 34738160
$arr[0,4] = 19041.6665
$arr[0,5] = -34735040
$arr[0,6] = -25281.6665
$ws.Range("H65:N65").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1645.875
$arr[0,1] = 1128.125
$arr[0,2] = 2163.625
$arr[0,3] = 3384.375
$arr[0,4] = 6490.875
$arr[0,5] = -849.375
$arr[0,6] = -11560.875
$ws.Range("H134:N134").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1298.5
$arr[0,1] = 768
$arr[0,2] = 2005.8334
$arr[0,3] = 2304
$arr[0,4] = 6017.5002
$arr[0,5] = -2192
$arr[0,6] = -6241.5002
$ws.Range("H5:N5").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3044.7144
$arr[0,1] = 2000
$arr[0,2] = 3218.8333
$arr[0,3] = 6000
$arr[0,4] = 9656.499899999999
$arr[0,5] = -5872
$arr[0,6] = -9912.499899999999
$ws.Range("H58:N58").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 550
$arr[0,1] = 375
$arr[0,2] = 666.6667
$arr[0,3] = 1125
$arr[0,4] = 2000.0001
$arr[0,5] = -314
$arr[0,6] = -3622.0001
$ws.Range("H68:N68").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 550
$arr[0,1] = 375
$arr[0,2] = 666.6667
$arr[0,3] = 3375
$arr[0,4] = 6000.0003
$arr[0,5] = 681
$arr[0,6] = -14112.0003
$ws.Range("H71:N71").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 560
$arr[0,1] = 433.33334
$arr[0,2] = 750
$arr[0,3] = 1300.00002
$arr[0,4] = 2250
$arr[0,5] = -52.00001999999995
$arr[0,6] = -4746
$ws.Range("H92:N92").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1077.8064
$arr[0,1] = 365.7143
$arr[0,2] = 1285.5
$arr[0,3] = 3291.4287
$arr[0,4] = 11569.5
$arr[0,5] = -841.4286999999999
$arr[0,6] = -16469.5
$ws.Range("H122:N122").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2688.25
$arr[0,1] = 30030
$arr[0,2] = 1987.1794
$arr[0,3] = 90090
$arr[0,4] = 5961.5382
$arr[0,5] = -85050
$arr[0,6] = -16041.5382
$ws.Range("H131:N131").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3888.889
$arr[0,1] = 4517.3335
$arr[0,2] = 2632
$arr[0,3] = 40656.0015
$arr[0,4] = 23688
$arr[0,5] = -38126.0015
$arr[0,6] = -28748
$ws.Range("H132:N132").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1298.5
$arr[0,1] = 768
$arr[0,2] = 2005.8334
$arr[0,3] = 6912
$arr[0,4] = 18052.5006
$arr[0,5] = -4377
$arr[0,6] = -23122.5006
$ws.Range("H135:N135").Value = $arr

$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1381
$arr[0,1] = 1088.75
$arr[0,2] = 2550
$arr[0,3] = 1088.75
$arr[0,4] = 2550
$arr[0,5] = -592.75
$arr[0,6] = -3542
$ws.Range("H97:N97").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6239.0586
$arr[0,1] = 8567
$arr[0,2] = 652
$arr[0,3] = 8567
$arr[0,4] = 652
$arr[0,5] = -6647
$arr[0,6] = -4492
$ws.Range("H107:N107").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2675.7222
$arr[0,1] = 2088.7778
$arr[0,2] = 3262.6667
$arr[0,3] = 6266.3334
$arr[0,4] = 9788.000100000001
$arr[0,5] = -3796.3334
$arr[0,6] = -14728.0001
$ws.Range("H126:N126").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5435.92
$arr[0,1] = 6541.6665
$arr[0,2] = 4415.231
$arr[0,3] = 19624.9995
$arr[0,4] = 13245.693
$arr[0,5] = -17094.9995
$arr[0,6] = -18305.693
$ws.Range("H132:N132").Value = $arr

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1779
$arr[0,1] = 1450
$arr[0,2] = 2272.5
$arr[0,3] = 1450
$arr[0,4] = 2272.5
$arr[0,5] = -1262
$arr[0,6] = -2648.5
$ws.Range("H46:N46").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 52635744
$arr[0,1] = 5241.6665
$arr[0,2] = 142859460
$arr[0,3] = 5241.6665
$arr[0,4] = 142859460
$arr[0,5] = -4700.6665
$arr[0,6] = -142860542
$ws.Range("H100:N100").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10423257
$arr[0,1] = 12823063
$arr[0,2] = 24097.223
$arr[0,3] = 38469189
$arr[0,4] = 72291.66900000001
$arr[0,5] = -38466659
$arr[0,6] = -77351.66900000001
$ws.Range("H132:N132").Value = $arr

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 15662.417
$arr[0,1] = 4500
$arr[0,2] = 19383.223
$arr[0,3] = 4500
$arr[0,4] = 19383.223
$arr[0,5] = -3980
$arr[0,6] = -20423.223
$ws.Range("H54:N54").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 18383760
$arr[0,1] = 24416178
$arr[0,2] = 717393.2
$arr[0,3] = 73248534
$arr[0,4] = 2152179.6
$arr[0,5] = -73245984
$arr[0,6] = -2157279.6
$ws.Range("H136:N136").Value = $arr
